$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.279.20"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "'3.111.33"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'575.06"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'178.07"
$ws.Range("E6").Value = "  +6.08%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'3.107.79"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'6.50"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").Value = "'0.152"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "'0.467"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'36.39"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "'3.631.99"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "'67.296.31"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'3.110.87"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").Value = "'16.55"
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").Value = "'485.75"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'7.72"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'83.66"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "'12.77"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'2.26"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").Value = "'10.37"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'7.90"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "'2.31"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").Value = "'2.61"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'28.12"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").Value = "'0.0₃0939"
$ws.Range("E34").Value = "  +3.56%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'47.56"
$ws.Range("E36").Value = "  +3.06%  "
$ws.Range("D37").Value = "'0.946"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").Value = "'0.316"
$ws.Range("E39").Value = "  +4.30%  "
$ws.Range("D40").Value = "'49.23"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("E44").Value = "  +8.44%  "
$ws.Range("D45").Value = "'2.783.26"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").Value = "'372.80"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'26.56"
$ws.Range("E49").Value = "  +8.00%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "'2.34"
$ws.Range("E51").Value = "  +8.66%  "
